# Auto update Excel log
# Appends new sensor-log rows to the PIR, Humidity and Proximity sheets,
# matching the rows added upstream (all dated 2026-01-30, hour "17:00").

$wb = $excel.ActiveWorkbook

function Append-Rows($SheetName, $StartRow, $Rows) {

    $ws = $wb.Worksheets.Item($SheetName)

    # Column A holds date-looking text ("2026-01-30") and column E on the
    # Humidity sheet holds percentage-looking text ("86.8%"). Both would be
    # auto-coerced into numeric/date values by plain .Value assignment, so
    # force those ranges to Text format first to keep them literal strings
    # (matching the rest of the column).
    $endRow = $StartRow + $Rows.Length - 1
    $ws.Range("A" + $StartRow + ":A" + $endRow).NumberFormat = "@"
    if ($SheetName -eq "Humidity") {
        $ws.Range("E" + $StartRow + ":E" + $endRow).NumberFormat = "@"
    }

    for ($i = 0; $i -lt $Rows.Length; $i++) {
        $r = $StartRow + $i
        $vals = $Rows[$i]
        for ($c = 0; $c -lt $vals.Length; $c++) {
            $ws.Cells.Item($r, $c + 1).Value = $vals[$c]
        }
    }
}

# --- PIR sheet: rows 402-414 (dimension A1:F401 -> A1:F414) ---
$pirRows = @(
    @("2026-01-30","17:56:18","17:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","17:56:20","17:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","17:56:23","17:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","17:56:28","17:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","17:56:33","17:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","17:56:38","17:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","17:56:43","17:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","17:56:48","17:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","17:56:53","17:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","17:56:58","17:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","17:57:03","17:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","17:57:08","17:00","Bathroom","No Motion","Inactive"),
    @("2026-01-30","17:57:13","17:00","Bathroom","No Motion","Inactive")
)
Append-Rows "PIR" 402 $pirRows

# --- Humidity sheet: rows 271-276 (dimension A1:F270 -> A1:F276) ---
$humidityRows = @(
    @("2026-01-30","17:56:21","17:00","Bathroom","86.8%","Active"),
    @("2026-01-30","17:56:29","17:00","Bathroom","86.8%","Active"),
    @("2026-01-30","17:56:39","17:00","Bathroom","86.8%","Active"),
    @("2026-01-30","17:56:59","17:00","Bathroom","86.7%","Active"),
    @("2026-01-30","17:57:04","17:00","Bathroom","86.8%","Active"),
    @("2026-01-30","17:57:09","17:00","Bathroom","86.8%","Active")
)
Append-Rows "Humidity" 271 $humidityRows

# --- Proximity sheet: rows 74-83 (dimension A1:F73 -> A1:F83) ---
$proximityRows = @(
    @("2026-01-30","17:56:19","17:00","Bathroom Door","ENTER","User ENTERED Bathroom"),
    @("2026-01-30","17:56:21","17:00","Bathroom Door","EXIT","User EXITED Bathroom"),
    @("2026-01-30","17:56:26","17:00","Bathroom Door","ENTER","User ENTERED Bathroom"),
    @("2026-01-30","17:56:29","17:00","Bathroom Door","EXIT","User EXITED Bathroom"),
    @("2026-01-30","17:56:35","17:00","Bathroom Door","ENTER","User ENTERED Bathroom"),
    @("2026-01-30","17:56:40","17:00","Bathroom Door","EXIT","User EXITED Bathroom"),
    @("2026-01-30","17:56:45","17:00","Bathroom Door","ENTER","User ENTERED Bathroom"),
    @("2026-01-30","17:56:49","17:00","Bathroom Door","EXIT","User EXITED Bathroom"),
    @("2026-01-30","17:56:51","17:00","Bathroom Door","ENTER","User ENTERED Bathroom"),
    @("2026-01-30","17:57:11","17:00","Bathroom Door","EXIT","User EXITED Bathroom")
)
Append-Rows "Proximity" 74 $proximityRows
